$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows continuing the time series through 2021-07-25 (serial 44402),
# appended after the existing last row (301, serial 44375).
$data = @(
    @(44376,0,0,0),
    @(44377,0,0,0),
    @(44378,0,0,0),
    @(44379,0,0,0),
    @(44380,0,0,0),
    @(44381,0,0,0),
    @(44382,0,0,0),
    @(44383,0,0,0),
    @(44384,0,0,0),
    @(44385,0,0,0),
    @(44386,3,3,19.79414093428345),
    @(44387,0,3,19.79414093428345),
    @(44388,1,4,26.39218791237794),
    @(44389,0,4,26.39218791237794),
    @(44390,0,4,26.39218791237794),
    @(44391,0,4,26.39218791237794),
    @(44392,0,4,26.39218791237794),
    @(44393,0,1,6.598046978094485),
    @(44394,0,1,6.598046978094485),
    @(44395,1,1,6.598046978094485),
    @(44396,0,1,6.598046978094485),
    @(44397,0,1,6.598046978094485),
    @(44398,0,1,6.598046978094485),
    @(44399,0,1,6.598046978094485),
    @(44400,3,4,26.39218791237794),
    @(44401,0,4,26.39218791237794),
    @(44402,0,3,19.79414093428345)
)

$startRow = 302
$lastExistingRow = 301

# Template cell for column A carries style s="2" (date number format, border,
# bold, centered). Copy format+value, then overwrite the value per row so the
# style (and not just the number format) is reproduced exactly.
$srcA = $ws.Range($ws.Cells.Item($lastExistingRow, 1), $ws.Cells.Item($lastExistingRow, 1))

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]

    $dstA = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 1))
    $srcA.Copy($dstA)
    $ws.Cells.Item($r, 1).Value = $vals[0]

    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}
